# Auto-generated script: extend Gold Coast stats sheet from column KN to KS (round 9 data)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Apply the same (default/general) style used by column KM to the existing KN column
#    and to the five brand-new columns KO:KR (matches the "s=1"-equivalent styling the
#    author applied when they extended the sheet). The final new column (KS) is left
#    unstyled, matching the source data exactly.
$styleRange = $ws.Range($ws.Cells.Item(1, 300), $ws.Cells.Item(102, 304))
$styleRange.Style = $ws.Cells.Item(1, 299).Style

# 2. Write the new match data for columns KO (301), KP (302), KQ (303), KR (304), KS (305)
#    Column KN (300) keeps its existing value -- only its style changes (handled above).
$ws.Cells.Item(1, 301).Value2 = 11009
$ws.Cells.Item(1, 302).Value2 = 11019
$ws.Cells.Item(1, 303).Value2 = 11028
$ws.Cells.Item(1, 304).Value2 = 11038
$ws.Cells.Item(1, 305).Value2 = 11045
$ws.Cells.Item(2, 301).Value2 = 2023
$ws.Cells.Item(2, 302).Value2 = 2023
$ws.Cells.Item(2, 303).Value2 = 2023
$ws.Cells.Item(2, 304).Value2 = 2023
$ws.Cells.Item(2, 305).Value2 = 2023
$ws.Cells.Item(3, 301).Value2 = 5
$ws.Cells.Item(3, 302).Value2 = 6
$ws.Cells.Item(3, 303).Value2 = 7
$ws.Cells.Item(3, 304).Value2 = 8
$ws.Cells.Item(3, 305).Value2 = 9
$ws.Cells.Item(4, 301).Value2 = 0
$ws.Cells.Item(4, 302).Value2 = 1
$ws.Cells.Item(4, 303).Value2 = 0
$ws.Cells.Item(4, 304).Value2 = 1
$ws.Cells.Item(4, 305).Value2 = 0
$ws.Cells.Item(5, 301).Value2 = 0
$ws.Cells.Item(5, 302).Value2 = 0
$ws.Cells.Item(5, 303).Value2 = 0
$ws.Cells.Item(5, 304).Value2 = 0
$ws.Cells.Item(5, 305).Value2 = 0
$ws.Cells.Item(6, 301).Value2 = 109
$ws.Cells.Item(6, 302).Value2 = 57
$ws.Cells.Item(6, 303).Value2 = 112
$ws.Cells.Item(6, 304).Value2 = 45
$ws.Cells.Item(6, 305).Value2 = 120
$ws.Cells.Item(7, 301).Value2 = 56
$ws.Cells.Item(7, 302).Value2 = 110
$ws.Cells.Item(7, 303).Value2 = 75
$ws.Cells.Item(7, 304).Value2 = 79
$ws.Cells.Item(7, 305).Value2 = 52
$ws.Cells.Item(8, 301).Value2 = 53
$ws.Cells.Item(8, 302).Value2 = -53
$ws.Cells.Item(8, 303).Value2 = 37
$ws.Cells.Item(8, 304).Value2 = -34
$ws.Cells.Item(8, 305).Value2 = 68
$ws.Cells.Item(9, 301).Value2 = 1
$ws.Cells.Item(9, 302).Value2 = 0
$ws.Cells.Item(9, 303).Value2 = 1
$ws.Cells.Item(9, 304).Value2 = 0
$ws.Cells.Item(9, 305).Value2 = 1
$ws.Cells.Item(10, 301).Value2 = 10
$ws.Cells.Item(10, 302).Value2 = 16
$ws.Cells.Item(10, 303).Value2 = 17
$ws.Cells.Item(10, 304).Value2 = 2
$ws.Cells.Item(10, 305).Value2 = 12
$ws.Cells.Item(11, 301).Value2 = 234
$ws.Cells.Item(11, 302).Value2 = 192
$ws.Cells.Item(11, 303).Value2 = 195
$ws.Cells.Item(11, 304).Value2 = 167
$ws.Cells.Item(11, 305).Value2 = 224
$ws.Cells.Item(12, 301).Value2 = 132
$ws.Cells.Item(12, 302).Value2 = 141
$ws.Cells.Item(12, 303).Value2 = 127
$ws.Cells.Item(12, 304).Value2 = 137
$ws.Cells.Item(12, 305).Value2 = 158
$ws.Cells.Item(13, 301).Value2 = 366
$ws.Cells.Item(13, 302).Value2 = 333
$ws.Cells.Item(13, 303).Value2 = 322
$ws.Cells.Item(13, 304).Value2 = 304
$ws.Cells.Item(13, 305).Value2 = 382
$ws.Cells.Item(14, 301).Value2 = 1.77
$ws.Cells.Item(14, 302).Value2 = 1.36
$ws.Cells.Item(14, 303).Value2 = 1.54
$ws.Cells.Item(14, 304).Value2 = 1.22
$ws.Cells.Item(14, 305).Value2 = 1.42
$ws.Cells.Item(15, 301).Value2 = 110
$ws.Cells.Item(15, 302).Value2 = 85
$ws.Cells.Item(15, 303).Value2 = 79
$ws.Cells.Item(15, 304).Value2 = 60
$ws.Cells.Item(15, 305).Value2 = 63
$ws.Cells.Item(16, 301).Value2 = 61
$ws.Cells.Item(16, 302).Value2 = 62
$ws.Cells.Item(16, 303).Value2 = 79
$ws.Cells.Item(16, 304).Value2 = 63
$ws.Cells.Item(16, 305).Value2 = 78
$ws.Cells.Item(17, 301).Value2 = 39
$ws.Cells.Item(17, 302).Value2 = 37
$ws.Cells.Item(17, 303).Value2 = 63
$ws.Cells.Item(17, 304).Value2 = 46
$ws.Cells.Item(17, 305).Value2 = 45
$ws.Cells.Item(18, 301).Value2 = 27
$ws.Cells.Item(18, 302).Value2 = 19
$ws.Cells.Item(18, 303).Value2 = 30
$ws.Cells.Item(18, 304).Value2 = 17
$ws.Cells.Item(18, 305).Value2 = 25
$ws.Cells.Item(19, 301).Value2 = 20
$ws.Cells.Item(19, 302).Value2 = 15
$ws.Cells.Item(19, 303).Value2 = 18
$ws.Cells.Item(19, 304).Value2 = 21
$ws.Cells.Item(19, 305).Value2 = 14
$ws.Cells.Item(20, 301).Value2 = 16
$ws.Cells.Item(20, 302).Value2 = 8
$ws.Cells.Item(20, 303).Value2 = 17
$ws.Cells.Item(20, 304).Value2 = 6
$ws.Cells.Item(20, 305).Value2 = 17
$ws.Cells.Item(21, 301).Value2 = 10
$ws.Cells.Item(21, 302).Value2 = 6
$ws.Cells.Item(21, 303).Value2 = 13
$ws.Cells.Item(21, 304).Value2 = 4
$ws.Cells.Item(21, 305).Value2 = 10
$ws.Cells.Item(22, 301).Value2 = 8
$ws.Cells.Item(22, 302).Value2 = 8
$ws.Cells.Item(22, 303).Value2 = 10
$ws.Cells.Item(22, 304).Value2 = 7
$ws.Cells.Item(22, 305).Value2 = 14
$ws.Cells.Item(23, 301).Value2 = 5
$ws.Cells.Item(23, 302).Value2 = 1
$ws.Cells.Item(23, 303).Value2 = 0
$ws.Cells.Item(23, 304).Value2 = 2
$ws.Cells.Item(23, 305).Value2 = 4
$ws.Cells.Item(24, 301).Value2 = 29
$ws.Cells.Item(24, 302).Value2 = 17
$ws.Cells.Item(24, 303).Value2 = 27
$ws.Cells.Item(24, 304).Value2 = 15
$ws.Cells.Item(24, 305).Value2 = 35
$ws.Cells.Item(25, 301).Value2 = 55.2
$ws.Cells.Item(25, 302).Value2 = 47.1
$ws.Cells.Item(25, 303).Value2 = 63
$ws.Cells.Item(25, 304).Value2 = 40
$ws.Cells.Item(25, 305).Value2 = 48.6
$ws.Cells.Item(26, 301).Value2 = 22.88
$ws.Cells.Item(26, 302).Value2 = 41.62
$ws.Cells.Item(26, 303).Value2 = 18.94
$ws.Cells.Item(26, 304).Value2 = 50.67
$ws.Cells.Item(26, 305).Value2 = 22.47
$ws.Cells.Item(27, 301).Value2 = 12.62
$ws.Cells.Item(27, 302).Value2 = 19.59
$ws.Cells.Item(27, 303).Value2 = 11.93
$ws.Cells.Item(27, 304).Value2 = 20.27
$ws.Cells.Item(27, 305).Value2 = 10.91
$ws.Cells.Item(28, 301).Value2 = 38
$ws.Cells.Item(28, 302).Value2 = 30
$ws.Cells.Item(28, 303).Value2 = 48
$ws.Cells.Item(28, 304).Value2 = 34
$ws.Cells.Item(28, 305).Value2 = 40
$ws.Cells.Item(29, 301).Value2 = 48
$ws.Cells.Item(29, 302).Value2 = 59
$ws.Cells.Item(29, 303).Value2 = 47
$ws.Cells.Item(29, 304).Value2 = 71
$ws.Cells.Item(29, 305).Value2 = 51
$ws.Cells.Item(30, 301).Value2 = 36
$ws.Cells.Item(30, 302).Value2 = 40
$ws.Cells.Item(30, 303).Value2 = 22
$ws.Cells.Item(30, 304).Value2 = 54
$ws.Cells.Item(30, 305).Value2 = 36
$ws.Cells.Item(31, 301).Value2 = 65
$ws.Cells.Item(31, 302).Value2 = 49
$ws.Cells.Item(31, 303).Value2 = 67
$ws.Cells.Item(31, 304).Value2 = 46
$ws.Cells.Item(31, 305).Value2 = 67
$ws.Cells.Item(32, 301).Value2 = 2.24
$ws.Cells.Item(32, 302).Value2 = 2.88
$ws.Cells.Item(32, 303).Value2 = 2.48
$ws.Cells.Item(32, 304).Value2 = 3.07
$ws.Cells.Item(32, 305).Value2 = 1.91
$ws.Cells.Item(33, 301).Value2 = 4.06
$ws.Cells.Item(33, 302).Value2 = 6.12
$ws.Cells.Item(33, 303).Value2 = 3.94
$ws.Cells.Item(33, 304).Value2 = 7.67
$ws.Cells.Item(33, 305).Value2 = 3.94
$ws.Cells.Item(34, 301).Value2 = 36.9
$ws.Cells.Item(34, 302).Value2 = 32.7
$ws.Cells.Item(34, 303).Value2 = 40.3
$ws.Cells.Item(34, 304).Value2 = 28.3
$ws.Cells.Item(34, 305).Value2 = 46.3
$ws.Cells.Item(35, 301).Value2 = 24.6
$ws.Cells.Item(35, 302).Value2 = 16.3
$ws.Cells.Item(35, 303).Value2 = 25.4
$ws.Cells.Item(35, 304).Value2 = 13
$ws.Cells.Item(35, 305).Value2 = 25.4
$ws.Cells.Item(36, 301).Value2 = 190
$ws.Cells.Item(36, 302).Value2 = 189.6
$ws.Cells.Item(36, 303).Value2 = 189.2
$ws.Cells.Item(36, 304).Value2 = 189.5
$ws.Cells.Item(36, 305).Value2 = 189.3
$ws.Cells.Item(37, 301).Value2 = 87.3
$ws.Cells.Item(37, 302).Value2 = 87.09999999999999
$ws.Cells.Item(37, 303).Value2 = 87.40000000000001
$ws.Cells.Item(37, 304).Value2 = 86.59999999999999
$ws.Cells.Item(37, 305).Value2 = 86.5
$ws.Cells.Item(38, 301).Value2 = 24.8
$ws.Cells.Item(38, 302).Value2 = 23.91
$ws.Cells.Item(38, 303).Value2 = 24.91
$ws.Cells.Item(38, 304).Value2 = 24.66
$ws.Cells.Item(38, 305).Value2 = 24
$ws.Cells.Item(39, 301).Value2 = 75.59999999999999
$ws.Cells.Item(39, 302).Value2 = 74.8
$ws.Cells.Item(39, 303).Value2 = 87.7
$ws.Cells.Item(39, 304).Value2 = 85.3
$ws.Cells.Item(39, 305).Value2 = 72.2
$ws.Cells.Item(40, 301).Value2 = 10
$ws.Cells.Item(40, 302).Value2 = 9
$ws.Cells.Item(40, 303).Value2 = 7
$ws.Cells.Item(40, 304).Value2 = 7
$ws.Cells.Item(40, 305).Value2 = 10
$ws.Cells.Item(41, 301).Value2 = 5
$ws.Cells.Item(41, 302).Value2 = 6
$ws.Cells.Item(41, 303).Value2 = 5
$ws.Cells.Item(41, 304).Value2 = 6
$ws.Cells.Item(41, 305).Value2 = 5
$ws.Cells.Item(42, 301).Value2 = 5
$ws.Cells.Item(42, 302).Value2 = 5
$ws.Cells.Item(42, 303).Value2 = 7
$ws.Cells.Item(42, 304).Value2 = 7
$ws.Cells.Item(42, 305).Value2 = 5
$ws.Cells.Item(43, 301).Value2 = 3
$ws.Cells.Item(43, 302).Value2 = 3
$ws.Cells.Item(43, 303).Value2 = 4
$ws.Cells.Item(43, 304).Value2 = 3
$ws.Cells.Item(43, 305).Value2 = 3
$ws.Cells.Item(44, 301).Value2 = 135
$ws.Cells.Item(44, 302).Value2 = 125
$ws.Cells.Item(44, 303).Value2 = 148
$ws.Cells.Item(44, 304).Value2 = 116
$ws.Cells.Item(44, 305).Value2 = 155
$ws.Cells.Item(45, 301).Value2 = 231
$ws.Cells.Item(45, 302).Value2 = 205
$ws.Cells.Item(45, 303).Value2 = 176
$ws.Cells.Item(45, 304).Value2 = 174
$ws.Cells.Item(45, 305).Value2 = 211
$ws.Cells.Item(46, 301).Value2 = 279
$ws.Cells.Item(46, 302).Value2 = 236
$ws.Cells.Item(46, 303).Value2 = 230
$ws.Cells.Item(46, 304).Value2 = 220
$ws.Cells.Item(46, 305).Value2 = 280
$ws.Cells.Item(47, 301).Value2 = 76.2
$ws.Cells.Item(47, 302).Value2 = 70.90000000000001
$ws.Cells.Item(47, 303).Value2 = 71.40000000000001
$ws.Cells.Item(47, 304).Value2 = 72.40000000000001
$ws.Cells.Item(47, 305).Value2 = 73.3
$ws.Cells.Item(48, 301).Value2 = 48
$ws.Cells.Item(48, 302).Value2 = 59
$ws.Cells.Item(48, 303).Value2 = 47
$ws.Cells.Item(48, 304).Value2 = 71
$ws.Cells.Item(48, 305).Value2 = 51
$ws.Cells.Item(49, 301).Value2 = 6
$ws.Cells.Item(49, 302).Value2 = 9
$ws.Cells.Item(49, 303).Value2 = 11
$ws.Cells.Item(49, 304).Value2 = 8
$ws.Cells.Item(49, 305).Value2 = 9
$ws.Cells.Item(50, 301).Value2 = 8
$ws.Cells.Item(50, 302).Value2 = 5
$ws.Cells.Item(50, 303).Value2 = 12
$ws.Cells.Item(50, 304).Value2 = 6
$ws.Cells.Item(50, 305).Value2 = 8
$ws.Cells.Item(51, 301).Value2 = 38
$ws.Cells.Item(51, 302).Value2 = 30
$ws.Cells.Item(51, 303).Value2 = 48
$ws.Cells.Item(51, 304).Value2 = 34
$ws.Cells.Item(51, 305).Value2 = 40
$ws.Cells.Item(52, 301).Value2 = 36
$ws.Cells.Item(52, 302).Value2 = 40
$ws.Cells.Item(52, 303).Value2 = 22
$ws.Cells.Item(52, 304).Value2 = 54
$ws.Cells.Item(52, 305).Value2 = 36
$ws.Cells.Item(53, 301).Value2 = 44
$ws.Cells.Item(53, 302).Value2 = 43
$ws.Cells.Item(53, 303).Value2 = 39
$ws.Cells.Item(53, 304).Value2 = 37
$ws.Cells.Item(53, 305).Value2 = 47
$ws.Cells.Item(54, 301).Value2 = 5
$ws.Cells.Item(54, 302).Value2 = 3
$ws.Cells.Item(54, 303).Value2 = 6
$ws.Cells.Item(54, 304).Value2 = 5
$ws.Cells.Item(54, 305).Value2 = 8
$ws.Cells.Item(55, 301).Value2 = 10
$ws.Cells.Item(55, 302).Value2 = 6
$ws.Cells.Item(55, 303).Value2 = 13
$ws.Cells.Item(55, 304).Value2 = 4
$ws.Cells.Item(55, 305).Value2 = 10
$ws.Cells.Item(56, 301).Value2 = 62.5
$ws.Cells.Item(56, 302).Value2 = 75
$ws.Cells.Item(56, 303).Value2 = 76.5
$ws.Cells.Item(56, 304).Value2 = 66.7
$ws.Cells.Item(56, 305).Value2 = 58.8
$ws.Cells.Item(57, 301).Value2 = 169
$ws.Cells.Item(57, 302).Value2 = 205
$ws.Cells.Item(57, 303).Value2 = 186
$ws.Cells.Item(57, 304).Value2 = 269
$ws.Cells.Item(57, 305).Value2 = 238
$ws.Cells.Item(58, 301).Value2 = 101
$ws.Cells.Item(58, 302).Value2 = 162
$ws.Cells.Item(58, 303).Value2 = 90
$ws.Cells.Item(58, 304).Value2 = 108
$ws.Cells.Item(58, 305).Value2 = 154
$ws.Cells.Item(59, 301).Value2 = 270
$ws.Cells.Item(59, 302).Value2 = 367
$ws.Cells.Item(59, 303).Value2 = 276
$ws.Cells.Item(59, 304).Value2 = 377
$ws.Cells.Item(59, 305).Value2 = 392
$ws.Cells.Item(60, 301).Value2 = 1.67
$ws.Cells.Item(60, 302).Value2 = 1.27
$ws.Cells.Item(60, 303).Value2 = 2.07
$ws.Cells.Item(60, 304).Value2 = 2.49
$ws.Cells.Item(60, 305).Value2 = 1.55
$ws.Cells.Item(61, 301).Value2 = 64
$ws.Cells.Item(61, 302).Value2 = 103
$ws.Cells.Item(61, 303).Value2 = 78
$ws.Cells.Item(61, 304).Value2 = 148
$ws.Cells.Item(61, 305).Value2 = 73
$ws.Cells.Item(62, 301).Value2 = 44
$ws.Cells.Item(62, 302).Value2 = 73
$ws.Cells.Item(62, 303).Value2 = 67
$ws.Cells.Item(62, 304).Value2 = 60
$ws.Cells.Item(62, 305).Value2 = 56
$ws.Cells.Item(63, 301).Value2 = 30
$ws.Cells.Item(63, 302).Value2 = 38
$ws.Cells.Item(63, 303).Value2 = 18
$ws.Cells.Item(63, 304).Value2 = 37
$ws.Cells.Item(63, 305).Value2 = 25
$ws.Cells.Item(64, 301).Value2 = 20
$ws.Cells.Item(64, 302).Value2 = 15
$ws.Cells.Item(64, 303).Value2 = 18
$ws.Cells.Item(64, 304).Value2 = 21
$ws.Cells.Item(64, 305).Value2 = 14
$ws.Cells.Item(65, 301).Value2 = 27
$ws.Cells.Item(65, 302).Value2 = 18
$ws.Cells.Item(65, 303).Value2 = 30
$ws.Cells.Item(65, 304).Value2 = 18
$ws.Cells.Item(65, 305).Value2 = 25
$ws.Cells.Item(66, 301).Value2 = 8
$ws.Cells.Item(66, 302).Value2 = 17
$ws.Cells.Item(66, 303).Value2 = 12
$ws.Cells.Item(66, 304).Value2 = 10
$ws.Cells.Item(66, 305).Value2 = 7
$ws.Cells.Item(67, 301).Value2 = 4
$ws.Cells.Item(67, 302).Value2 = 12
$ws.Cells.Item(67, 303).Value2 = 6
$ws.Cells.Item(67, 304).Value2 = 3
$ws.Cells.Item(67, 305).Value2 = 4
$ws.Cells.Item(68, 301).Value2 = 8
$ws.Cells.Item(68, 302).Value2 = 8
$ws.Cells.Item(68, 303).Value2 = 3
$ws.Cells.Item(68, 304).Value2 = 14
$ws.Cells.Item(68, 305).Value2 = 9
$ws.Cells.Item(69, 301).Value2 = 0
$ws.Cells.Item(69, 302).Value2 = 0
$ws.Cells.Item(69, 303).Value2 = 0
$ws.Cells.Item(69, 304).Value2 = 5
$ws.Cells.Item(69, 305).Value2 = 1
$ws.Cells.Item(70, 301).Value2 = 16
$ws.Cells.Item(70, 302).Value2 = 25
$ws.Cells.Item(70, 303).Value2 = 15
$ws.Cells.Item(70, 304).Value2 = 29
$ws.Cells.Item(70, 305).Value2 = 17
$ws.Cells.Item(71, 301).Value2 = 50
$ws.Cells.Item(71, 302).Value2 = 68
$ws.Cells.Item(71, 303).Value2 = 80
$ws.Cells.Item(71, 304).Value2 = 34.5
$ws.Cells.Item(71, 305).Value2 = 41.2
$ws.Cells.Item(72, 301).Value2 = 33.75
$ws.Cells.Item(72, 302).Value2 = 21.59
$ws.Cells.Item(72, 303).Value2 = 23
$ws.Cells.Item(72, 304).Value2 = 37.7
$ws.Cells.Item(72, 305).Value2 = 56
$ws.Cells.Item(73, 301).Value2 = 16.88
$ws.Cells.Item(73, 302).Value2 = 14.68
$ws.Cells.Item(73, 303).Value2 = 18.4
$ws.Cells.Item(73, 304).Value2 = 13
$ws.Cells.Item(73, 305).Value2 = 23.06
$ws.Cells.Item(74, 301).Value2 = 28
$ws.Cells.Item(74, 302).Value2 = 33
$ws.Cells.Item(74, 303).Value2 = 31
$ws.Cells.Item(74, 304).Value2 = 37
$ws.Cells.Item(74, 305).Value2 = 31
$ws.Cells.Item(75, 301).Value2 = 71
$ws.Cells.Item(75, 302).Value2 = 59
$ws.Cells.Item(75, 303).Value2 = 63
$ws.Cells.Item(75, 304).Value2 = 60
$ws.Cells.Item(75, 305).Value2 = 75
$ws.Cells.Item(76, 301).Value2 = 48
$ws.Cells.Item(76, 302).Value2 = 40
$ws.Cells.Item(76, 303).Value2 = 49
$ws.Cells.Item(76, 304).Value2 = 39
$ws.Cells.Item(76, 305).Value2 = 50
$ws.Cells.Item(77, 301).Value2 = 45
$ws.Cells.Item(77, 302).Value2 = 57
$ws.Cells.Item(77, 303).Value2 = 35
$ws.Cells.Item(77, 304).Value2 = 64
$ws.Cells.Item(77, 305).Value2 = 45
$ws.Cells.Item(78, 301).Value2 = 2.81
$ws.Cells.Item(78, 302).Value2 = 2.28
$ws.Cells.Item(78, 303).Value2 = 2.33
$ws.Cells.Item(78, 304).Value2 = 2.21
$ws.Cells.Item(78, 305).Value2 = 2.65
$ws.Cells.Item(79, 301).Value2 = 5.62
$ws.Cells.Item(79, 302).Value2 = 3.35
$ws.Cells.Item(79, 303).Value2 = 2.92
$ws.Cells.Item(79, 304).Value2 = 6.4
$ws.Cells.Item(79, 305).Value2 = 6.43
$ws.Cells.Item(80, 301).Value2 = 35.6
$ws.Cells.Item(80, 302).Value2 = 43.9
$ws.Cells.Item(80, 303).Value2 = 42.9
$ws.Cells.Item(80, 304).Value2 = 37.5
$ws.Cells.Item(80, 305).Value2 = 35.6
$ws.Cells.Item(81, 301).Value2 = 17.8
$ws.Cells.Item(81, 302).Value2 = 29.8
$ws.Cells.Item(81, 303).Value2 = 34.3
$ws.Cells.Item(81, 304).Value2 = 15.6
$ws.Cells.Item(81, 305).Value2 = 15.6
$ws.Cells.Item(82, 301).Value2 = 188.3
$ws.Cells.Item(82, 302).Value2 = 187
$ws.Cells.Item(82, 303).Value2 = 190
$ws.Cells.Item(82, 304).Value2 = 188.3
$ws.Cells.Item(82, 305).Value2 = 186.7
$ws.Cells.Item(83, 301).Value2 = 86.09999999999999
$ws.Cells.Item(83, 302).Value2 = 86.7
$ws.Cells.Item(83, 303).Value2 = 87.40000000000001
$ws.Cells.Item(83, 304).Value2 = 86.8
$ws.Cells.Item(83, 305).Value2 = 83.90000000000001
$ws.Cells.Item(84, 301).Value2 = 24
$ws.Cells.Item(84, 302).Value2 = 25.66
$ws.Cells.Item(84, 303).Value2 = 25.74
$ws.Cells.Item(84, 304).Value2 = 26.41
$ws.Cells.Item(84, 305).Value2 = 24.24
$ws.Cells.Item(85, 301).Value2 = 67.90000000000001
$ws.Cells.Item(85, 302).Value2 = 108.9
$ws.Cells.Item(85, 303).Value2 = 99.40000000000001
$ws.Cells.Item(85, 304).Value2 = 121.1
$ws.Cells.Item(85, 305).Value2 = 78.2
$ws.Cells.Item(86, 301).Value2 = 12
$ws.Cells.Item(86, 302).Value2 = 2
$ws.Cells.Item(86, 303).Value2 = 8
$ws.Cells.Item(86, 304).Value2 = 6
$ws.Cells.Item(86, 305).Value2 = 11
$ws.Cells.Item(87, 301).Value2 = 5
$ws.Cells.Item(87, 302).Value2 = 11
$ws.Cells.Item(87, 303).Value2 = 6
$ws.Cells.Item(87, 304).Value2 = 2
$ws.Cells.Item(87, 305).Value2 = 3
$ws.Cells.Item(88, 301).Value2 = 2
$ws.Cells.Item(88, 302).Value2 = 4
$ws.Cells.Item(88, 303).Value2 = 3
$ws.Cells.Item(88, 304).Value2 = 5
$ws.Cells.Item(88, 305).Value2 = 7
$ws.Cells.Item(89, 301).Value2 = 4
$ws.Cells.Item(89, 302).Value2 = 6
$ws.Cells.Item(89, 303).Value2 = 6
$ws.Cells.Item(89, 304).Value2 = 10
$ws.Cells.Item(89, 305).Value2 = 2
$ws.Cells.Item(90, 301).Value2 = 127
$ws.Cells.Item(90, 302).Value2 = 118
$ws.Cells.Item(90, 303).Value2 = 111
$ws.Cells.Item(90, 304).Value2 = 132
$ws.Cells.Item(90, 305).Value2 = 137
$ws.Cells.Item(91, 301).Value2 = 138
$ws.Cells.Item(91, 302).Value2 = 246
$ws.Cells.Item(91, 303).Value2 = 149
$ws.Cells.Item(91, 304).Value2 = 238
$ws.Cells.Item(91, 305).Value2 = 232
$ws.Cells.Item(92, 301).Value2 = 189
$ws.Cells.Item(92, 302).Value2 = 272
$ws.Cells.Item(92, 303).Value2 = 190
$ws.Cells.Item(92, 304).Value2 = 273
$ws.Cells.Item(92, 305).Value2 = 281
$ws.Cells.Item(93, 301).Value2 = 70
$ws.Cells.Item(93, 302).Value2 = 74.09999999999999
$ws.Cells.Item(93, 303).Value2 = 68.8
$ws.Cells.Item(93, 304).Value2 = 72.40000000000001
$ws.Cells.Item(93, 305).Value2 = 71.7
$ws.Cells.Item(94, 301).Value2 = 71
$ws.Cells.Item(94, 302).Value2 = 59
$ws.Cells.Item(94, 303).Value2 = 63
$ws.Cells.Item(94, 304).Value2 = 60
$ws.Cells.Item(94, 305).Value2 = 75
$ws.Cells.Item(95, 301).Value2 = 7
$ws.Cells.Item(95, 302).Value2 = 7
$ws.Cells.Item(95, 303).Value2 = 7
$ws.Cells.Item(95, 304).Value2 = 15
$ws.Cells.Item(95, 305).Value2 = 6
$ws.Cells.Item(96, 301).Value2 = 4
$ws.Cells.Item(96, 302).Value2 = 16
$ws.Cells.Item(96, 303).Value2 = 6
$ws.Cells.Item(96, 304).Value2 = 14
$ws.Cells.Item(96, 305).Value2 = 6
$ws.Cells.Item(97, 301).Value2 = 28
$ws.Cells.Item(97, 302).Value2 = 33
$ws.Cells.Item(97, 303).Value2 = 31
$ws.Cells.Item(97, 304).Value2 = 37
$ws.Cells.Item(97, 305).Value2 = 31
$ws.Cells.Item(98, 301).Value2 = 48
$ws.Cells.Item(98, 302).Value2 = 40
$ws.Cells.Item(98, 303).Value2 = 49
$ws.Cells.Item(98, 304).Value2 = 39
$ws.Cells.Item(98, 305).Value2 = 50
$ws.Cells.Item(99, 301).Value2 = 38
$ws.Cells.Item(99, 302).Value2 = 42
$ws.Cells.Item(99, 303).Value2 = 45
$ws.Cells.Item(99, 304).Value2 = 42
$ws.Cells.Item(99, 305).Value2 = 59
$ws.Cells.Item(100, 301).Value2 = 2
$ws.Cells.Item(100, 302).Value2 = 7
$ws.Cells.Item(100, 303).Value2 = 5
$ws.Cells.Item(100, 304).Value2 = 2
$ws.Cells.Item(100, 305).Value2 = 7
$ws.Cells.Item(101, 301).Value2 = 4
$ws.Cells.Item(101, 302).Value2 = 12
$ws.Cells.Item(101, 303).Value2 = 6
$ws.Cells.Item(101, 304).Value2 = 3
$ws.Cells.Item(101, 305).Value2 = 4
$ws.Cells.Item(102, 301).Value2 = 50
$ws.Cells.Item(102, 302).Value2 = 70.59999999999999
$ws.Cells.Item(102, 303).Value2 = 50
$ws.Cells.Item(102, 304).Value2 = 30
$ws.Cells.Item(102, 305).Value2 = 57.1
